$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Cells.Item(5, 2).Value = 317
$ws.Cells.Item(5, 3).Value = "cleveland"
$ws.Cells.Item(5, 4).Value = "cleveland@gmail.com"
$ws.Cells.Item(5, 5).Value = "nasa corporation"
$ws.Cells.Item(5, 6).Value = "'55667788"
$ws.Cells.Item(5, 7).Value = "L2_scheduled"

# Row 6
$ws.Cells.Item(6, 2).Value = 320
$ws.Cells.Item(6, 3).Value = "harsh"
$ws.Cells.Item(6, 4).Value = "harsh@mai.com"
$ws.Cells.Item(6, 5).Value = "wipro"
$ws.Cells.Item(6, 6).Value = "'123456789"
$ws.Cells.Item(6, 7).Value = "L2_TBS"
